$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.280823128560371
$ws.Range("C2").Value = 0.07915547753766106
$ws.Range("D2").Value = 0.5346783271571667
$ws.Range("E2").Value = 0.1768911601100633
$ws.Range("G2").Value = 0.002554510007074008
$ws.Range("J2").Value = 0.06741007856803183
$ws.Range("K2").Value = 0.7037056246823283
$ws.Range("L2").Value = 0.4189190924949031
$ws.Range("N2").Value = 2.777633187389654
$ws.Range("O2").Value = 7.141456727666366

$ws.Range("B3").Value = 1.244599554801709
$ws.Range("C3").Value = 0.07781329502025613
$ws.Range("D3").Value = 0.5322067500233061
$ws.Range("E3").Value = 0.1770125259674025
$ws.Range("G3").Value = 0.002557469389171832
$ws.Range("J3").Value = 0.06770158817154215
$ws.Range("K3").Value = 0.6700692321717838
$ws.Range("L3").Value = 0.415077664569381
$ws.Range("N3").Value = 2.800698875694088
$ws.Range("O3").Value = 7.167007916831977

$ws.Range("B4").Value = 1.222927364422105
$ws.Range("C4").Value = 0.07697751722614044
$ws.Range("D4").Value = 0.5309111381388192
$ws.Range("E4").Value = 0.1771483094187616
$ws.Range("G4").Value = 0.002559384908082829
$ws.Range("J4").Value = 0.06789281316686768
$ws.Range("K4").Value = 0.6497159731215163
$ws.Range("L4").Value = 0.4128930193174654
$ws.Range("N4").Value = 2.815603981677544
$ws.Range("O4").Value = 7.185539548630601

$ws.Range("B5").Value = 1.214239590865645
$ws.Range("C5").Value = 0.07663400386281438
$ws.Range("D5").Value = 0.5304390995114403
$ws.Range("E5").Value = 0.1772190809108309
$ws.Range("G5").Value = 0.002560190327128201
$ws.Range("J5").Value = 0.06797382260062967
$ws.Range("K5").Value = 0.6414976889402624
$ws.Range("L5").Value = 0.4120466366570454
$ws.Range("N5").Value = 2.821864743567115
$ws.Range("O5").Value = 7.193806831519794

$ws.Range("B6").Value = 1.212805699173771
$ws.Range("C6").Value = 0.07657678717634298
$ws.Range("D6").Value = 0.5303640997564116
$ws.Range("E6").Value = 0.1772317658971563
$ws.Range("G6").Value = 0.00256032556806391
$ws.Range("J6").Value = 0.06798746059859706
$ws.Range("K6").Value = 0.6401376438103057
$ws.Range("L6").Value = 0.4119087491970248
$ws.Range("N6").Value = 2.822915622912493
$ws.Range("O6").Value = 7.195222838466947

$ws.Range("B7").Value = 1.222809614843698
$ws.Range("C7").Value = 0.07697289632658766
$ws.Range("D7").Value = 0.5309045454224446
$ws.Range("E7").Value = 0.1771492013131333
$ws.Range("G7").Value = 0.002559395669705167
$ws.Range("J7").Value = 0.06789389319360684
$ws.Range("K7").Value = 0.6496048307654974
$ws.Range("L7").Value = 0.412881426880773
$ws.Range("N7").Value = 2.815687660041601
$ws.Range("O7").Value = 7.185648146346864

$ws.Range("B8").Value = 1.268215491574637
$ws.Range("C8").Value = 0.07869512001816048
$ws.Range("D8").Value = 0.5337801344661841
$ws.Range("E8").Value = 0.1769203112567261
$ws.Range("G8").Value = 0.002555510014262186
$ws.Range("J8").Value = 0.06750805529095594
$ws.Range("K8").Value = 0.6920458930116524
$ws.Range("L8").Value = 0.417558524332847
$ws.Range("N8").Value = 2.785432129849941
$ws.Range("O8").Value = 7.14967702720287

$ws.Range("B9").Value = 1.361748659136481
$ws.Range("C9").Value = 0.08197956880184165
$ws.Range("D9").Value = 0.5411756048577843
$ws.Range("E9").Value = 0.1769562348779452
$ws.Range("G9").Value = 0.002548667994948915
$ws.Range("J9").Value = 0.06684822552582936
$ws.Range("K9").Value = 0.7776332707523181
$ws.Range("L9").Value = 0.4281064666407559
$ws.Range("N9").Value = 2.731990557941373
$ws.Range("O9").Value = 7.101675256421998

$ws.Range("B10").Value = 1.433182139582868
$ws.Range("C10").Value = 0.08433598229149197
$ws.Range("D10").Value = 0.5476748809821004
$ws.Range("E10").Value = 0.1772765678741308
$ws.Range("G10").Value = 0.002544110550130263
$ws.Range("J10").Value = 0.0664220496761132
$ws.Range("K10").Value = 0.8419379558439459
$ws.Range("L10").Value = 0.4366902471509917
$ws.Range("N10").Value = 2.696311506239201
$ws.Range("O10").Value = 7.080125088032901

$ws.Range("B11").Value = 1.466263989204492
$ws.Range("C11").Value = 0.08539566530215836
$ws.Range("D11").Value = 0.5508620354911784
$ws.Range("E11").Value = 0.1774857988537271
$ws.Range("G11").Value = 0.002542138164688815
$ws.Range("J11").Value = 0.06624080945931432
$ws.Range("K11").Value = 0.8714981208953532
$ws.Range("L11").Value = 0.4407754525002758
$ws.Range("N11").Value = 2.680857270482985
$ws.Range("O11").Value = 7.07329552736212

$ws.Range("B12").Value = 1.478874980998
$ws.Range("C12").Value = 0.08579517117670576
$ws.Range("D12").Value = 0.5521019831228386
$ws.Range("E12").Value = 0.1775741317597657
$ws.Range("G12").Value = 0.002541405694860604
$ws.Range("J12").Value = 0.06617398795004092
$ws.Range("K12").Value = 0.8827356276290459
$ws.Range("L12").Value = 0.4423482512890473
$ws.Range("N12").Value = 2.675116742142055
$ws.Range("O12").Value = 7.07113653979934

$ws.Range("B13").Value = 1.476155270854974
$ws.Range("C13").Value = 0.08570920945496141
$ws.Range("D13").Value = 0.5518334705806183
$ws.Range("E13").Value = 0.1775547033806397
$ws.Range("G13").Value = 0.002541562804685905
$ws.Range("J13").Value = 0.06618829874007837
$ws.Range("K13").Value = 0.8803134936578942
$ws.Range("L13").Value = 0.4420083744968935
$ws.Range("N13").Value = 2.676348101404781
$ws.Range("O13").Value = 7.071582522017934

$ws.Range("B14").Value = 1.467299831420348
$ws.Range("C14").Value = 0.08542856848249158
$ws.Range("D14").Value = 0.5509633852097124
$ws.Range("E14").Value = 0.1774928838209462
$ws.Range("G14").Value = 0.002542077615089562
$ws.Range("J14").Value = 0.06623527576749133
$ws.Range("K14").Value = 0.8724217642527208
$ws.Range("L14").Value = 0.4409043308217377
$ws.Range("N14").Value = 2.680382755906621
$ws.Range("O14").Value = 7.073109346319484

$ws.Range("B15").Value = 1.461886486919013
$ws.Range("C15").Value = 0.08525643658121851
$ws.Range("D15").Value = 0.5504347322396939
$ws.Range("E15").Value = 0.1774562019451231
$ws.Range("G15").Value = 0.002542394828663525
$ws.Range("J15").Value = 0.06626428612502089
$ws.Range("K15").Value = 0.8675935293797465
$ws.Range("L15").Value = 0.4402314311488453
$ws.Range("N15").Value = 2.682868639702708
$ws.Range("O15").Value = 7.074100195050846

$ws.Range("B16").Value = 1.43103190204954
$ws.Range("C16").Value = 0.08426648193771769
$ws.Range("D16").Value = 0.5474712230676744
$ws.Range("E16").Value = 0.1772641694691366
$ws.Range("G16").Value = 0.002544241473076685
$ws.Range("J16").Value = 0.06643414778303391
$ws.Range("K16").Value = 0.8400122752093182
$ws.Range("L16").Value = 0.4364268891292369
$ws.Range("N16").Value = 2.697337100375591
$ws.Range("O16").Value = 7.080631220554892

$ws.Range("B17").Value = 1.412253299262431
$ws.Range("C17").Value = 0.08365602934676986
$ws.Range("D17").Value = 0.5457121841104993
$ws.Range("E17").Value = 0.1771626081024529
$ws.Range("G17").Value = 0.002545400103627804
$ws.Range("J17").Value = 0.06654158289979684
$ws.Range("K17").Value = 0.823170490019919
$ws.Range("L17").Value = 0.434139044434346
$ws.Range("N17").Value = 2.706411884374841
$ws.Range("O17").Value = 7.085399149043496

$ws.Range("B18").Value = 1.401507580880718
$ws.Range("C18").Value = 0.08330375966323089
$ws.Range("D18").Value = 0.5447221437961929
$ws.Range("E18").Value = 0.177110172387593
$ws.Range("G18").Value = 0.002546076010589266
$ws.Range("J18").Value = 0.06660456581554985
$ws.Range("K18").Value = 0.8135125286214304
$ws.Range("L18").Value = 0.4328401240637163
$ws.Range("N18").Value = 2.711704531947419
$ws.Range("O18").Value = 7.088421459382744

$ws.Range("B19").Value = 1.397878773987827
$ws.Range("C19").Value = 0.0831842893887611
$ws.Range("D19").Value = 0.544390665900039
$ws.Range("E19").Value = 0.1770934465195708
$ws.Range("G19").Value = 0.002546306493797673
$ws.Range("J19").Value = 0.06662609515310347
$ws.Range("K19").Value = 0.810247503547572
$ws.Range("L19").Value = 0.4324032532409348
$ws.Range("N19").Value = 2.713509081147862
$ws.Range("O19").Value = 7.089492851179045

$ws.Range("B20").Value = 1.414246601885992
$ws.Range("C20").Value = 0.08372113252218583
$ws.Range("D20").Value = 0.5458971905233483
$ws.Range("E20").Value = 0.1771728008318867
$ws.Range("G20").Value = 0.002545275783268821
$ws.Range("J20").Value = 0.06653002322563761
$ws.Range("K20").Value = 0.8249603304579978
$ws.Range("L20").Value = 0.4343808321562221
$ws.Range("N20").Value = 2.705438295340603
$ws.Range("O20").Value = 7.084862626723066

$ws.Range("B21").Value = 1.469898625058704
$ws.Range("C21").Value = 0.08551104773801654
$ws.Range("D21").Value = 0.5512180545863146
$ws.Range("E21").Value = 0.1775107949537293
$ws.Range("G21").Value = 0.002541926011719868
$ws.Range("J21").Value = 0.06622142839318634
$ws.Range("K21").Value = 0.8747385742782399
$ws.Range("L21").Value = 0.4412279154239513
$ws.Range("N21").Value = 2.679194649622747
$ws.Range("O21").Value = 7.072649289565646

$ws.Range("B22").Value = 1.506757472535242
$ws.Range("C22").Value = 0.08667052133588271
$ws.Range("D22").Value = 0.5548880823384934
$ws.Range("E22").Value = 0.1777847308807559
$ws.Range("G22").Value = 0.002539820819654527
$ws.Range("J22").Value = 0.06603029253446735
$ws.Range("K22").Value = 0.9075260882479483
$ws.Range("L22").Value = 0.445853323855161
$ws.Range("N22").Value = 2.662693733667233
$ws.Range("O22").Value = 7.067157144563168

$ws.Range("B23").Value = 1.487040882931637
$ws.Range("C23").Value = 0.08605263778968464
$ws.Range("D23").Value = 0.5529117416496092
$ws.Range("E23").Value = 0.1776336834161398
$ws.Range("G23").Value = 0.002540936729314138
$ws.Range("J23").Value = 0.06613134205453974
$ws.Range("K23").Value = 0.8900036565377434
$ws.Range("L23").Value = 0.4433709307039067
$ws.Range("N23").Value = 2.671441025796643
$ws.Range("O23").Value = 7.069860703071924

$ws.Range("B24").Value = 1.413345272878388
$ws.Range("C24").Value = 0.08369170346366417
$ws.Range("D24").Value = 0.5458134829103614
$ws.Range("E24").Value = 0.1771681741530422
$ws.Range("G24").Value = 0.002545331957889612
$ws.Range("J24").Value = 0.0665352455627195
$ws.Range("K24").Value = 0.8241510669135153
$ws.Range("L24").Value = 0.4342714688594072
$ws.Range("N24").Value = 2.705878219965246
$ws.Range("O24").Value = 7.085104312606319

$ws.Range("B25").Value = 1.335967004937231
$ws.Range("C25").Value = 0.08110098461802551
$ws.Range("D25").Value = 0.5389874071873209
$ws.Range("E25").Value = 0.1768947830263805
$ws.Range("G25").Value = 0.002550436172445216
$ws.Range("J25").Value = 0.06701640661235864
$ws.Range("K25").Value = 0.7542284984106118
$ws.Range("L25").Value = 0.425106131636511
$ws.Range("N25").Value = 2.745818036997871
$ws.Range("O25").Value = 7.11225058358383
Write-Output "Updated 240 cells"
